$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACLgrantHistTable")

# Update the descriptive text cells to reflect the new fiscal-year range
# (FY 2011-2016 -> FY 2012-2016), keeping all other content/formatting intact.
$ws.Range("A3").Value = "This table shows the grant awards and award dollars ACL made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the ACL page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars ACL made for FY 2012-2016."
